$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily NAV rows appended to the bottom of the data table (rows 735-754),
# one row per trading day from 2024-09-02 through 2024-09-27.
# Columns: A=Date(text) C..H=constituent prices/basket value I=daily return J=NAV
$csv = @"
2024-09-02,1885.400024414062,683.5999755859375,73.80999755859375,296.8999938964844,1303.849975585938,28487.57955932617,0,483.8294342156934
2024-09-03,1901.949951171875,689.4000244140625,74.47000122070312,297.1499938964844,1320.25,28715.86001586914,0.008013332830455756,487.7065205052348
2024-09-04,1901.300048828125,688.9500122070312,74.16000366210938,298.9500122070312,1327.099975585938,28721.03076171875,0.0001800658537390795,487.7943397962237
2024-09-05,1879.449951171875,687.5,76,290.6000061035156,1312.349975585938,28602.74978637695,-0.004118270556621158,485.7854707289544
2024-09-06,1872.349975585938,673.5499877929688,74.72000122070312,283.6000061035156,1289.699951171875,28191.60983276367,-0.01437414083205038,478.8027219585325
2024-09-09,1892.400024414062,664.1500244140625,74.33999633789062,281.5499877929688,1237.150024414062,28036.46997070312,-0.005503050836077006,476.1678462391426
2024-09-10,1922.449951171875,664.5999755859375,78.05000305175781,285.75,1250.300048828125,28561.49984741211,0.01872667555001102,485.0848870030106
2024-09-11,1957.599975585938,689.75,81.94999694824219,288.0499877929688,1237.699951171875,29297.64938354492,0.0257741904334731,497.587557257026
2024-09-12,1996.400024414062,729.1500244140625,81.69999694824219,291.7000122070312,1237.300048828125,29912.80038452148,0.02099659917843318,508.0352037529274
2024-09-13,1988.050048828125,713.7000122070312,83.11000061035156,289.9500122070312,1241.5,29812.18057250977,-0.003363771051799782,506.3262896412481
2024-09-16,1989.900024414062,714.2000122070312,84.69999694824219,290.3999938964844,1226.599975585938,29926.49987792969,0.003834650911961043,508.2678742095708
2024-09-17,2006.550048828125,731.0999755859375,82,284.2999877929688,1193.800048828125,29823.24993896484,-0.003450117433913109,506.5142903556625
2024-09-18,1987.800048828125,743.25,80.81999969482422,282.8500061035156,1166.400024414062,29685.31034851074,-0.004625236710834788,504.171541865347
2024-09-19,1998.599975585938,735.9500122070312,80.97000122070312,272.7000122070312,1121.300048828125,29455.8603515625,-0.0077294120982553,500.274592250057
2024-09-20,2048.10009765625,746.5,83.44999694824219,277.3500061035156,1149.400024414062,30118.95037841797,0.02251131078642199,511.5364290747486
2024-09-23,2082.39990234375,773.9500122070312,82.88999938964844,286.2999877929688,1162.75,30664.31942749023,0.01810717313253569,520.7989077596042
2024-09-24,2068.14990234375,781.8499755859375,83.79000091552734,291.7999877929688,1141.199951171875,30770.6690826416,0.003468188994144963,522.6051367996588
2024-09-25,2061.60009765625,775.8499755859375,82.95999908447266,289.8500061035156,1118.449951171875,30542.33015441895,-0.007420668286718119,518.7270574345335
2024-09-26,2022.050048828125,772.0999755859375,81.83999633789062,290.5,1123.650024414062,30238.56976318359,-0.009945553914831306,513.5680295177366
2024-09-27,2031.300048828125,761.75,81.08999633789062,293.4500122070312,1130.25,30174.3701171875,-0.002123104581297322,512.4776708814596
"@

$lines = $csv -split "`n"
$startRow = 735
$endRow = $startRow + $lines.Length - 1

# Force column A to text first so the "yyyy-mm-dd" date strings are stored
# as literal text (matching the rest of the sheet) instead of being
# auto-converted into date serial numbers.
$dateRange = $ws.Range("A$startRow`:A$endRow")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $lines.Length; $i++) {
    $line = $lines[$i].Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line.Split(",")
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value  = $parts[0]              # A: Date
    $ws.Cells.Item($r, 3).Value  = [double]$parts[1]      # C
    $ws.Cells.Item($r, 4).Value  = [double]$parts[2]      # D
    $ws.Cells.Item($r, 5).Value  = [double]$parts[3]      # E
    $ws.Cells.Item($r, 6).Value  = [double]$parts[4]      # F
    $ws.Cells.Item($r, 7).Value  = [double]$parts[5]      # G
    $ws.Cells.Item($r, 8).Value  = [double]$parts[6]      # H
    $ws.Cells.Item($r, 9).Value  = [double]$parts[7]      # I
    $ws.Cells.Item($r, 10).Value = [double]$parts[8]      # J
}

# Drop the temporary text format again so the new cells carry no explicit
# style, same as the existing data rows above them.
$dateRange.ClearFormats()

Write-Output "Inserted $($lines.Length) rows ($startRow..$endRow)"
